# Extra support for temp dirs for optimize tests
#
# The "optimize" worksheet contained a key/value settings table with a
# row for "domain_type" = "segmentation". That row is removed here,
# which shifts the remaining rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimize")
$ws.Activate()

# Find and delete the row containing the "domain_type" key in column A.
$found = $ws.Columns.Item(1).Find("domain_type")
if ($found -ne $null) {
    $found.EntireRow.Select()
    $found.EntireRow.Delete()
}
else {
    # Fallback: scan column A manually for the "domain_type" row.
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value -eq "domain_type") {
            $ws.Rows.Item($r).Select()
            $ws.Rows.Item($r).Delete()
            break
        }
    }
}

# Select the row that now occupies row 3 (previously row 4), matching
# the resulting workbook state where row 3 is selected on the optimize tab.
$ws.Rows.Item(3).Select()
